$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 188
$ws.Range("F5").Value = 1035
$ws.Range("F7").Value = 2669
$ws.Range("F9").Value = 1303
$ws.Range("F10").Value = 934
$ws.Range("F11").Value = 629
$ws.Range("F12").Value = 943
$ws.Range("F13").Value = 1187
$ws.Range("F15").Value = 119
$ws.Range("F17").Value = 794
$ws.Range("F19").Value = 531
$ws.Range("F20").Value = 1138
$ws.Range("F22").Value = 645
$ws.Range("F24").Value = 231
$ws.Range("F26").Value = 314
$ws.Range("F27").Value = 697
$ws.Range("F28").Value = 585
$ws.Range("F29").Value = 5732
$ws.Range("F30").Value = 499
$ws.Range("F31").Value = 43
$ws.Range("F34").Value = 181
$ws.Range("F35").Value = 1646
$ws.Range("F37").Value = 99
$ws.Range("F38").Value = 448
$ws.Range("F41").Value = 153
$ws.Range("F42").Value = 14
$ws.Range("F45").Value = 146
$ws.Range("F47").Value = 121

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 88
$ws.Range("F9").Value = 51
$ws.Range("F12").Value = 196
$ws.Range("F13").Value = 4412
$ws.Range("F14").Value = 40
$ws.Range("F17").Value = 40
$ws.Range("F18").Value = 216

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 748

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 748
$ws.Range("F5").Value = 88
$ws.Range("F6").Value = 1035
$ws.Range("F7").Value = 2670
$ws.Range("F9").Value = 1303
$ws.Range("F10").Value = 934
$ws.Range("F11").Value = 629
$ws.Range("F12").Value = 943
$ws.Range("F13").Value = 1187
$ws.Range("F16").Value = 119
$ws.Range("F19").Value = 794
$ws.Range("F21").Value = 531
$ws.Range("F22").Value = 1138
$ws.Range("F24").Value = 51
$ws.Range("F25").Value = 645
$ws.Range("F27").Value = 231
$ws.Range("F29").Value = 314
$ws.Range("F30").Value = 585
$ws.Range("F31").Value = 5732
$ws.Range("F32").Value = 196
$ws.Range("F33").Value = 499
$ws.Range("F34").Value = 43
$ws.Range("F36").Value = 181
$ws.Range("F37").Value = 1646
$ws.Range("F39").Value = 448
$ws.Range("F40").Value = 40
$ws.Range("F41").Value = 40
$ws.Range("F43").Value = 40
$ws.Range("F44").Value = 14
$ws.Range("F46").Value = 146
$ws.Range("F48").Value = 121

Write-Output "Applied 67 cell updates"
